$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Test Sheet 1" ---
$ws1 = $wb.Worksheets.Item("Test Sheet 1")
$ws1.Range("C10").Value = 12.1
$ws1.Range("C11").Select()

# --- Sheet 2: "Test Sheet 2" ---
$ws2 = $wb.Worksheets.Item("Test Sheet 2")
$ws2.Range("D9").Select()

# Restore Sheet 1 as the active/selected tab (matches original tabSelected state)
$ws1.Activate()
